$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10; existing rows 10-23 shift down to 11-24.
$ws.Rows("10:10").Insert()

# Populate the new row 10 with fresh data (values that were previously
# constant across the whole column block are copied from a neighboring row).
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Terminal La Palmera de La Serena"
$ws.Range("C10").Value = "Coquimbo"
$ws.Range("D10").Value = 44803
$ws.Range("E10").Value = 4
$ws.Range("F10").Value = 100112026
$ws.Range("G10").Value = "Haba"
$ws.Range("H10").Value = "Sin especificar"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 600
$ws.Range("K10").Value = 9500
$ws.Range("L10").Value = 10000
$ws.Range("M10").Value = 9750
$ws.Range("N10").Value = "$/saco 25 kilos"
$ws.Range("O10").Value = "Provincia del Elquí"
$ws.Range("P10").Value = 390
$ws.Range("Q10").Value = 25
$ws.Range("R10").Value = "Hortaliza"
